## ---------------------------------------------------------------------
## Add the 2022-Q4 quarter: a new summary row on "总计" and a brand new
## "2022-Q4" worksheet (holding the per-fund detail), inserted right
## after "总计" and before the existing "2022-Q3" sheet.
## ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

## -----------------------------------------------------------------
## 1) "总计" sheet: insert a new "2022-Q4" row above the existing
##    "2022-Q3" row (row 2), pushing everything else down by one.
## -----------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing data rows (2..7) down to rows (3..8), carrying the
# cell formatting (column A keeps its bold/bordered style) along.
$total.Range("A2:D7").Copy()
$total.Range("A3:D8").PasteSpecial()

# The freshly-uncovered last row (row 8) sometimes doesn't inherit the
# style from the paste above, so stamp it explicitly from its neighbour.
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)

$total.Application.CutCopyMode = $false

# Re-write the running index column for every row that moved down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6

# Fill in the brand-new 2022-Q4 summary row.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 26
$total.Range("D2").Value = 2.37

## -----------------------------------------------------------------
## 2) Create the new "2022-Q4" worksheet right after "总计" and
##    before "2022-Q3".
## -----------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Match the header-row / index-column look used by every other quarter
# sheet: bold text, thin box border, centered & top-aligned.
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$rows = @(
    @(0, "090001", "大成价值增长混合", "12.94", "66.15", "2.79", "0.3610", 10),
    @(1, "000574", "宝盈新价值灵活配置混合A", "7.64", "88.96", "4.17", "0.3186", 10),
    @(2, "003715", "宝盈消费主题灵活配置混合", "5.61", "83.13", "5.08", "0.2850", 5),
    @(3, "001915", "宝盈医疗健康沪港深股票", "6.60", "94.14", "4.14", "0.2732", 8),
    @(4, "000339", "长城医疗保健混合A", "7.61", "91.86", "2.57", "0.1956", 10),
    @(5, "090020", "大成健康产业混合A", "2.49", "92.90", "6.05", "0.1506", 6),
    @(6, "011673", "长城医药科技六个月持有期混合型证券投资基金A", "5.97", "91.30", "2.52", "0.1504", 10),
    @(7, "012045", "大成医药健康股票A", "2.04", "92.92", "6.24", "0.1273", 6),
    @(8, "007574", "宝盈新价值灵活配置混合C", "3.00", "88.96", "4.17", "0.1251", 10),
    @(9, "005347", "诺德量化优选6个月持有期混合", "2.06", "91.14", "4.07", "0.0838", 9),
    @(10, "006267", "诺德量化核心灵活配置混合A", "1.09", "91.67", "4.80", "0.0523", 5),
    @(11, "014020", "诺德量化先锋一年持有期混合A", "1.27", "72.32", "3.92", "0.0498", 9),
    @(12, "005293", "诺德新旺灵活配置混合", "0.55", "92.81", "8.01", "0.0441", 2),
    @(13, "010799", "长城优选稳进六个月持有期混合A", "1.96", "32.14", "1.17", "0.0229", 7),
    @(14, "006881", "华宝大健康混合", "0.79", "87.95", "2.86", "0.0226", 9),
    @(15, "006268", "诺德量化核心灵活配置混合C", "0.40", "91.67", "4.80", "0.0192", 5),
    @(16, "014021", "诺德量化先锋一年持有期混合C", "0.46", "72.32", "3.92", "0.0180", 9),
    @(17, "012046", "大成医药健康股票C", "0.27", "92.92", "6.24", "0.0168", 6),
    @(18, "010857", "宝盈祥乐一年持有期混合型证券投资基金A", "1.00", "33.11", "1.58", "0.0158", 4),
    @(19, "011674", "长城医药科技六个月持有期混合型证券投资基金C", "0.54", "91.30", "2.52", "0.0136", 10),
    @(20, "008324", "宝盈祥利稳健配置混合A", "0.53", "31.08", "1.86", "0.0099", 2),
    @(21, "008325", "宝盈祥利稳健配置混合C", "0.29", "31.08", "1.86", "0.0054", 2),
    @(22, "015562", "长城医疗保健混合C", "0.18", "91.86", "2.57", "0.0046", 10),
    @(23, "010858", "宝盈祥乐一年持有期混合型证券投资基金C", "0.06", "33.11", "1.58", "0.0009", 4),
    @(24, "016060", "大成健康产业混合C", "0.01", "92.90", "6.05", "0.0006", 6),
    @(25, "010800", "长城优选稳进六个月持有期混合C", "0.04", "32.14", "1.17", "0.0005", 7)
)

$rowIndex = 2
foreach ($r in $rows) {
    $idxCell = $q4.Cells.Item($rowIndex, 1)
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Value = $r[0]

    # Fund code must stay text (leading zeroes, e.g. "090001").
    $q4.Cells.Item($rowIndex, 2).NumberFormat = "@"
    $q4.Cells.Item($rowIndex, 2).Value = $r[1]

    $q4.Cells.Item($rowIndex, 3).Value = $r[2]

    # Columns D, E, F, G hold numeric-looking text ("12.94", not 12.94),
    # matching the source data's inlineStr cell type.
    $q4.Cells.Item($rowIndex, 4).NumberFormat = "@"
    $q4.Cells.Item($rowIndex, 4).Value = $r[3]

    $q4.Cells.Item($rowIndex, 5).NumberFormat = "@"
    $q4.Cells.Item($rowIndex, 5).Value = $r[4]

    $q4.Cells.Item($rowIndex, 6).NumberFormat = "@"
    $q4.Cells.Item($rowIndex, 6).Value = $r[5]

    $q4.Cells.Item($rowIndex, 7).NumberFormat = "@"
    $q4.Cells.Item($rowIndex, 7).Value = $r[6]

    $q4.Cells.Item($rowIndex, 8).Value = $r[7]

    $rowIndex++
}
